# Updates the cryptos table (column D = Price, column E = Volume(1h))
# Commit: Updated cryptos list on Sun Oct 22 07:37:56 UTC 2023 with GitHub Actions
#
# Note: D-column values are leading-apostrophe prefixed so Excel keeps them
# as text (as in the source file) rather than auto-converting the
# numeric-looking strings (e.g. "215.73") into actual numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.973.72"
$ws.Range("E2").Value = '  +1.08%  '

$ws.Range("D3").Value = "'1.639.83"
$ws.Range("E3").Value = '  +2.33%  '

$ws.Range("E4").Value = '  +0.37%  '

$ws.Range("D5").Value = "'215.73"
$ws.Range("E5").Value = '  +1.54%  '

$ws.Range("D6").Value = "'0.520"
$ws.Range("E6").Value = '  +0.61%  '

$ws.Range("E7").Value = '  +0.35%  '

$ws.Range("D8").Value = "'28.96"
$ws.Range("E8").Value = '  +3.79%  '

$ws.Range("E9").Value = '  +2.93%  '

$ws.Range("D10").Value = "'0.0611"
$ws.Range("E10").Value = '  +1.39%  '

$ws.Range("E11").Value = '  +0.79%  '

$ws.Range("D12").Value = "'1.878.63"
$ws.Range("E12").Value = '  +2.52%  '

$ws.Range("D13").Value = "'1.642.91"
$ws.Range("E13").Value = '  +2.86%  '

$ws.Range("D14").Value = "'0.570"
$ws.Range("E14").Value = '  +4.29%  '

$ws.Range("D15").Value = "'9.34"
$ws.Range("E15").Value = '  +19.95%  '

$ws.Range("D16").Value = "'3.91"
$ws.Range("E16").Value = '  +4.19%  '

$ws.Range("D17").Value = "'30.009.86"
$ws.Range("E17").Value = '  +1.12%  '

$ws.Range("D18").Value = "'64.67"
$ws.Range("E18").Value = '  +1.14%  '

$ws.Range("D19").Value = "'245.66"
$ws.Range("E19").Value = '  +1.30%  '

$ws.Range("D20").Value = "'0.0₃0706"
$ws.Range("E20").Value = '  +1.28%  '

$ws.Range("E21").Value = '  +0.20%  '

$ws.Range("D22").Value = "'10.02"
$ws.Range("E22").Value = '  +6.10%  '

$ws.Range("D23").Value = "'4.16"
$ws.Range("E23").Value = '  +3.59%  '

$ws.Range("D24").Value = "'2.16"
$ws.Range("E24").Value = '  +2.80%  '

$ws.Range("D25").Value = "'158.28"
$ws.Range("E25").Value = '  +1.83%  '

$ws.Range("D26").Value = "'15.68"
$ws.Range("E26").Value = '  +1.27%  '

$ws.Range("E27").Value = '  +2.10%  '

$ws.Range("D28").Value = "'6.67"
$ws.Range("E28").Value = '  +3.67%  '

$ws.Range("E29").Value = '  +0.38%  '

$ws.Range("D30").Value = "'0.0492"
$ws.Range("E30").Value = '  +2.15%  '

$ws.Range("E31").Value = '  +5.74%  '

$ws.Range("D32").Value = "'3.41"
$ws.Range("E32").Value = '  +5.41%  '

$ws.Range("D33").Value = "'3.19"
$ws.Range("E33").Value = '  +0.03%  '

$ws.Range("D34").Value = "'1.433.46"
$ws.Range("E34").Value = '  +0.53%  '

$ws.Range("D35").Value = "'1.67"
$ws.Range("E35").Value = '  +6.85%  '

$ws.Range("E36").Value = '  +1.58%  '

$ws.Range("E37").Value = '  -2.36%  '

$ws.Range("E38").Value = '  +0.35%  '

$ws.Range("E39").Value = '  +1.50%  '

$ws.Range("D40").Value = "'76.97"
$ws.Range("E40").Value = '  +15.93%  '

$ws.Range("D41").Value = "'0.558"
$ws.Range("E41").Value = '  +2.02%  '

$ws.Range("D42").Value = "'2.01"
$ws.Range("E42").Value = '  +1.99%  '

$ws.Range("D43").Value = "'0.838"
$ws.Range("E43").Value = '  +2.65%  '

$ws.Range("D44").Value = "'0.0496"
$ws.Range("E44").Value = '  -0.41%  '

$ws.Range("D45").Value = "'54.28"
$ws.Range("E45").Value = '  -7.17%  '

$ws.Range("E46").Value = '  +6.78%  '

$ws.Range("E47").Value = '  +0.29%  '

$ws.Range("D48").Value = "'5.40"
$ws.Range("E48").Value = '  +1.14%  '

$ws.Range("D49").Value = "'1.784.94"
$ws.Range("E49").Value = '  +2.29%  '

$ws.Range("E50").Value = '  +12.56%  '

$ws.Range("D51").Value = "'89.83"
$ws.Range("E51").Value = '  +3.47%  '
